$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 2 (old rows 2..59 shift down to 3..60).
$ws.Rows.Item(2).Insert()

# The inserted row inherits the header row's (row 1) bold/centered style by default.
# Clear that so it matches the plain style used by the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Column D holds dates rendered with the custom date number format used by every
# other data row in the sheet (style index 2 -> numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the brand-new row 2 with its data.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44631
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112012
$ws.Range("G2").Value = "Espinaca"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 9500
$ws.Range("M2").Value = 9273
$ws.Range("N2").Value = '$/cuna 10 kilos'
$ws.Range("O2").Value = "Provincia de Chacabuco"
$ws.Range("P2").Value = 927
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = "Hortaliza"
